$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'70.295.67"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = "'3.594.07"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'579.42"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.97%  '
$ws.Range('D6').Value = "'189.10"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.60%  '
$ws.Range('D7').Value = "'0.630"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.38%  '
$ws.Range('D8').Value = "'3.590.21"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.99%  '
$ws.Range('E10').Value = '  +3.13%  '
$ws.Range('D11').Value = "'0.662"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('D12').Value = "'55.89"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.21%  '
$ws.Range('D13').Value = "'0.0000311"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.73%  '
$ws.Range('D14').Value = "'9.67"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').Value = "'4.171.25"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = "'19.80"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').Value = "'3.592.44"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = "'70.218.11"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = "'12.66"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('E21').Value = '  -1.79%  '
$ws.Range('D22').Value = "'489.71"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = "'19.76"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').Value = "'4.88"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -9.41%  '
$ws.Range('D25').Value = "'97.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.46%  '
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('D27').Value = "'2.98"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.35%  '
$ws.Range('D28').Value = "'11.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.32%  '
$ws.Range('D29').Value = "'9.36"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').Value = "'32.27"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.53%  '
$ws.Range('D31').Value = "'7.61"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.62%  '
$ws.Range('D32').Value = "'12.23"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = "'0.118"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').Value = "'65.80"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = "'572.61"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.59%  '
$ws.Range('D36').Value = "'38.69"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.95%  '
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = "'0.0₃0807"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.82%  '
$ws.Range('D39').Value = "'0.398"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.90%  '
$ws.Range('D40').Value = "'2.97"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.27%  '
$ws.Range('D41').Value = "'3.23"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +15.75%  '
$ws.Range('D42').Value = "'3.49"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('E43').Value = '  -6.27%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = "'3.03"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.07%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = "'3.210.70"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.60%  '
$ws.Range('D46').Value = "'0.0445"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('D47').Value = "'3.49"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.53%  '
$ws.Range('D48').Value = "'9.57"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.30%  '
$ws.Range('D49').Value = "'0.138"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').Value = "'0.999"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('E51').Value = '  -4.35%  '
